$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Study 2")

# Insert 4 new columns (B:E) before the old "scale" column, shifting old B,C,D to F,G,H
$ws.Range("B1:E1").EntireColumn.Insert()

# Adjust column widths (closest achievable values in this environment)
$ws.Range("A:A").ColumnWidth = 17.5
$ws.Range("B:E").ColumnWidth = 14.5

$data = @(
  @('ess_variable_name', 'variable_label', 'ess4', 'ess7', 'ess8', 'scale', 'variable_label', 'response_scale'),
  @('idno', 'idno', 'yes', $null, $null, $null, 'Respondent''s identification number', $null),
  @('agea', 'age', 'yes', $null, $null, $null, 'Age of respondent, calculated', $null),
  @('gndr', 'gndr', 'yes', $null, $null, $null, 'Gender', $null),
  @('lrscale', 'lrscale', 'yes', $null, $null, $null, 'Placement on left right scale', $null),
  @('eisced', 'education_eisced', 'yes', $null, $null, $null, 'Highest level of education, ES - ISCED', $null),
  @('cntry', 'cntry', 'yes', $null, $null, $null, 'Country', $null),
  @('region', 'region', 'yes', $null, $null, $null, 'Region', $null),
  @('ipfrule', 'trad_1', 'yes', $null, $null, 'RWA', 'Important to do what is told and follow rules', $null),
  @('ipstrgv', 'trad_2', 'yes', $null, $null, 'RWA', 'Important that government is strong and ensures safety', $null),
  @('ipbhprp', 'trad_3', 'yes', $null, $null, 'RWA', 'Important to behave properly', $null),
  @('imptrad', 'trad_4', 'yes', $null, $null, 'RWA', 'Important to follow traditions and customs', $null),
  @('impsafe', 'trad_5', 'yes', $null, $null, 'RWA', 'Important to live in secure and safe surroundings', $null),
  @('hrshsnt', 'trad_6', 'yes', $null, $null, 'RWA', 'People who break the law much harsher sentences', $null),
  @('imdfetn', 'prj_immi_1', 'yes', $null, $null, 'Prejudice_immigrants', 'Allow many/few immigrants of different race/ethnic group from majority', $null),
  @('impcntr', 'prj_immi_2', 'yes', $null, $null, 'Prejudice_immigrants', 'Allow many/few immigrants from poorer countries outside Europe', $null),
  @('imbgeco', 'prj_immi_3', 'yes', $null, $null, 'Prejudice_immigrants', 'Immigration bad or good for country''s economy', $null),
  @('imueclt', 'prj_immi_4', 'yes', $null, $null, 'Prejudice_immigrants', 'Country''s cultural life undermined or enriched by immigrants', $null),
  @('imwbcnt', 'prj_immi_5', 'yes', $null, $null, 'Prejudice_immigrants', 'Immigrants make country worse or better place to live', $null),
  @('freehms', 'prj_gay_1', 'yes', $null, $null, 'Prejudice_gays', 'Gays and lesbians free to live life as they wish', $null),
  @('wmcpwrk', 'prj_wmn_1', 'yes', $null, $null, 'Prejudice_women', 'Women should be prepared to cut down on paid work for sake of family', $null),
  @('mnrgtjb', 'prj_wmn_2', 'yes', $null, $null, 'Prejudice_women', 'Men should have more right to job than women when jobs are scarce', $null),
  @('uentrjb', 'prj_unempl_1', 'yes', $null, $null, 'Prejudice_unemployed', 'Most unemployed people do not really try to find a job', $null),
  @('oafl20', 'prj_age20', 'yes', $null, $null, 'Prejudice_age', 'Overall how negative or positive feel towards people in their 20s', $null),
  @('oafl70', 'prj_age70', 'yes', $null, $null, 'Prejudice_age', 'Overall how negative or positive feel towards people over 70', $null),
  @('aljewlv', 'prj_jews', 'no', $null, $null, 'Prejudice_jews', 'Allow many or few Jewish people to come and live in country', $null),
  @('almuslv', 'prj_muslms', 'no', $null, $null, 'Prejudice_muslims', 'Allow many or few Muslims to come and live in country', $null),
  @('algyplv', 'prj_sinti', 'no', $null, $null, 'Prejudice_sinti_roma', 'Allow many or few Gypsies to come and live in country', $null),
  @('alpfpe', 'prj_prf_poor_1', 'no', $null, $null, 'Prejudice_poorer_countries', 'Allow professionals from [poor European country providing largest number of migrants]', $null),
  @('alpfpne', 'prj_prf_poor_2', 'no', $null, $null, 'Prejudice_poorer_countries', 'Allow professionals from [poor non-European country providing largest number of migrants]', $null),
  @('allbpe', 'prj_unsk_poor_1', 'no', $null, $null, 'Prejudice_poorer_countries', 'Allow unskilled labourers from [poor European country providing largest number of migrants]', $null),
  @('allbpne', 'prj_unsk_poor_2', 'no', $null, $null, 'Prejudice_poorer_countries', 'Allow unskilled labourers from [poor non-European country providing largest number of migrants]', $null),
  @('gvslvue', 'prj_unempl_2', $null, $null, $null, 'Prejudice_unemployed', $null, $null),
  @('dfincac', 'prj_unempl_3', $null, $null, $null, $null, $null, $null),
  @('bennent', 'prj_unempl_4', $null, $null, $null, $null, $null, $null),
  @('sblazy', 'prj_unempl_5', $null, $null, $null, $null, $null, $null),
  @('sbbsntx', 'prj_unempl_6', $null, $null, $null, $null, $null, $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowNum = $i + 1
  for ($j = 0; $j -lt $data[$i].Length; $j++) {
    $val = $data[$i][$j]
    if ($null -ne $val) {
      $ws.Cells.Item($rowNum, $j + 1).Value = $val
    }
  }
}

# Set the final selection to match the authored state
$ws.Activate()
$ws.Range("B37").Select()
